$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cells whose new text would otherwise be auto-converted to a number by Excel
# need their NumberFormat forced to text ("@") first, so the value round-trips
# as a string (matching the original inline-string cell type in the workbook).

$ws.Range("D2").Value = '61.524.85'
$ws.Range("E2").Value = '  -1.68%  '
$ws.Range("D3").Value = '2.994.43'
$ws.Range("E3").Value = '  -0.62%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '596.16'
$ws.Range("E5").Value = '  +1.96%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.84'
$ws.Range("E6").Value = '  -2.45%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").Value = '2.994.05'
$ws.Range("E9").Value = '  -0.56%  '
$ws.Range("E10").Value = '  -1.98%  '
$ws.Range("E11").Value = '  +2.05%  '
$ws.Range("E12").Value = '  +4.12%  '
$ws.Range("E13").Value = '  -0.68%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.27'
$ws.Range("E14").Value = '  -1.87%  '
$ws.Range("E15").Value = '  +2.10%  '
$ws.Range("D16").Value = '3.490.10'
$ws.Range("E16").Value = '  -0.72%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.03'
$ws.Range("E17").Value = '  +0.11%  '
$ws.Range("D18").Value = '61.489.32'
$ws.Range("E18").Value = '  -1.74%  '
$ws.Range("D19").Value = '2.993.99'
$ws.Range("E19").Value = '  -0.70%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '452.80'
$ws.Range("E20").Value = '  -2.80%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.96'
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("E22").Value = '  -0.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.34'
$ws.Range("E23").Value = '  -0.28%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '81.98'
$ws.Range("E24").Value = '  +2.18%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.24'
$ws.Range("E25").Value = '  -4.48%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.64'
$ws.Range("E26").Value = '  +3.63%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.97'
$ws.Range("E27").Value = '  -3.62%  '
$ws.Range("E28").Value = '  +0.02%  '
$ws.Range("E29").Value = '  +1.71%  '
$ws.Range("E30").Value = '  -0.03%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.23'
$ws.Range("E31").Value = '  +1.04%  '
$ws.Range("E32").Value = '  -2.64%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.46'
$ws.Range("E33").Value = '  -0.46%  '
$ws.Range("E34").Value = '  +2.30%  '
$ws.Range("D35").Value = '0.0₃0831'
$ws.Range("E35").Value = '  +4.37%  '
$ws.Range("E36").Value = '  -1.60%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.78'
$ws.Range("E37").Value = '  +0.49%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '9.27'
$ws.Range("E38").Value = '  +3.42%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '50.31'
$ws.Range("E39").Value = '  -0.17%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.07'
$ws.Range("E40").Value = '  -3.18%  '
$ws.Range("E41").Value = '  +9.59%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.89'
$ws.Range("E42").Value = '  -1.60%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '395.62'
$ws.Range("E43").Value = '  -6.27%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '39.62'
$ws.Range("E44").Value = '  +4.65%  '
$ws.Range("E45").Value = '  -0.25%  '
$ws.Range("E46").Value = '  -3.30%  '
$ws.Range("D47").Value = '2.713.58'
$ws.Range("E47").Value = '  -2.73%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '133.74'
$ws.Range("E48").Value = '  +3.31%  '
$ws.Range("E50").Value = '  -0.53%  '
$ws.Range("E51").Value = '  +1.44%  '
